$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 448 (shifts existing rows 448..481 down to 449..482,
# carrying their formatting, and grows the used range to R482).
$ws.Rows.Item(448).Insert()

# Populate the newly inserted row 448 with the new weekly price record.
$ws.Range("A448").Value = 8
$ws.Range("B448").Value = "Terminal La Palmera de La Serena"
$ws.Range("C448").Value = "Coquimbo"
$ws.Range("D448").Value = 45021
$ws.Range("E448").Value = 4
$ws.Range("F448").Value = 100112032
$ws.Range("G448").Value = "Zapallo italiano"
$ws.Range("H448").Value = "Sin especificar"
$ws.Range("I448").Value = "Primera"
$ws.Range("J448").Value = 480
$ws.Range("K448").Value = 7000
$ws.Range("L448").Value = 8000
$ws.Range("M448").Value = 7500
$ws.Range("N448").Value = "$/caja 70 unidades"
$ws.Range("O448").Value = "Provincia de Limarí"
$ws.Range("P448").Value = 107
$ws.Range("Q448").Value = 70
$ws.Range("R448").Value = "Hortaliza"
